# Refatorando artefatos seguindo a correção da ac5
# Adds a trailing period "." to the end of the N01, N02, N03, N04 and N05
# requirement paragraphs in the "Lista de Necessidades" list.

$d = $word.ActiveDocument

# --- N01: " Cadastro de fornecedores e seus produtos" -> add new run "."
$p = $d.Paragraphs(4)
$r = $p.Range
$r.End = $r.End - 1          # exclude the paragraph mark
$r.Collapse(0)                # wdCollapseEnd -> collapse to the end of the text
$r.InsertAfter(".")

# --- N02: ": Cadastro dos clientes da fábrica" -> add new run "."
$p = $d.Paragraphs(6)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(".")

# --- N03: " Gerenciamento de estoque " -> " Gerenciamento de estoque" + new run "."
$p = $d.Paragraphs(8)
$r = $p.Range
$r.End = $r.End - 1
$spaceRange = $d.Range($r.End - 1, $r.End)
$spaceRange.Text = ""         # drop the trailing space from the existing run

$p = $d.Paragraphs(8)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(".")

# --- N04: trailing " " run -> "." (same run, in place, no extra run added)
$p = $d.Paragraphs(10)
$r = $p.Range
$r.End = $r.End - 1
$lastCharRange = $d.Range($r.End - 1, $r.End)
$lastCharRange.Delete()
$lastCharRange.InsertAfter(".")

# --- N05: "Gestão de orçamentos e pedidos" -> add new run "."
$p = $d.Paragraphs(12)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(".")
